$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header row (B1:L1) and set A1 to the new title
$ws.Range("B1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# New column layout: A=Day, B=Time, C=Hours, D=Module Code, E=Module Title,
#                     F=Class Type, G=Lecturer, H=Group, I=Block, J=Room
# (columns K and L are no longer used)

$data = @(
    @("SUN", "9:30-11:30",  2,   "5CS022", "Human Computer Interaction",                 "Tutorial", "Mr. Dipesh Shrestha", "L5CG13",         "HCK", "TR-09  Chandragiri"),
    @("SUN", "13:00-15:00", 2,   "5CS020", "Distributed and Cloud Systems Programming",   "Tutorial", "Mr. Shishir Poudel",  "L5CG13",         "WLV", "SR-03 Wolves"),
    @("MON", "9:30-12:00",  2.5, "5CS022", "Human Computer Interaction",                 "Workshop", "Mr. Dipesh Shrestha", "L5CG13",         "HCK", "TR-09  Chandragiri"),
    @("MON", "13:00-15:30", 2.5, "5CS020", "Distributed and Cloud Systems Programming",   "Workshop", "Mr. Shishir Poudel",  "L5CG13",         "WLV", "Lab-02 Moseley"),
    @("TUE", "7:00-9:00",   2,   "5CS024", "Collaborative Development",                   "Lecture",  "Mr. Raj Shrestha",    "L5CG(12+13+14)", "WLV", "LT-03 Walsall"),
    @("WED", "7:00-9:00",   2,   "5CS022", "Human Computer Interaction",                 "Lecture",  "Mr. Ayush Shakya",    "L5CG(12+13+14)", "WLV", "LT-01 Wulfruna"),
    @("THU", "9:30-11:30",  2,   "5CS020", "Distributed and Cloud Systems Programming",   "Lecture",  "Mr. Sumanta Silwal",  "L5CG(12+13+14)", "WLV", "LT-01 Wulfruna"),
    @("THU", "13:00-15:00", 2,   "5CS024", "Collaborative Development",                   "Tutorial", "Mr. Udaya Kandel",    "L5CG13",         "WLV", "TR-01 Dudley"),
    @("FRI", "13:00-15:30", 2.5, "5CS024", "Collaborative Development",                   "Workshop", "Mr. Udaya Kandel",    "L5CG13",         "WLV", "SR-01 Bantok")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Remove the now-unused K and L columns' data (rows 2-10)
$ws.Range("K1:L10").ClearContents()
